# Insert a new paragraph ("哈哈哈，天气很不错呀！" with a gram-check mark
# wrapped around the middle "哈") right after the second occurrence of the
# paragraph "今天晴，今天早八，今天是开心的一天。" (i.e. immediately before the
# final "3月8号，踏青的好日子！" paragraph).

$d = $word.ActiveDocument

# Locate the target paragraph: the LAST paragraph whose text equals the
# target sentence (there are two identical paragraphs in the document; the
# new content goes after the second one).
$targetText = "今天晴，今天早八，今天是开心的一天。"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq $targetText) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $targetPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

# The freshly-created (empty) paragraph is now the next one.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range
$newRange.Collapse(0)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

$paragraphXml = '<w:p ' + $wNs + '>' `
    + '<w:pPr><w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' `
    + '<w:r>' + $rPr + '<w:t>哈哈</w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r>' + $rPr + '<w:t>哈</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r>' + $rPr + '<w:t>，天气很不错呀！</w:t></w:r>' `
    + '</w:p>'

$packageXml = '<?xml version="1.0" standalone="yes"?>' `
    + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $paragraphXml + '</w:body></w:document></pkg:xmlData>' `
    + '</pkg:part></pkg:package>'

$newRange.InsertXML($packageXml)

Write-Host "Inserted paragraph after paragraph index $targetIndex. Total paragraphs now:" $d.Paragraphs.Count
